$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Thbs1"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 31.645482
$ws.Cells.Item(2, 8).Value = 63.290964
$ws.Cells.Item(2, 9).Value = 0.03555980726701226
$ws.Cells.Item(2, 10).Value = 0.0244846141215985
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.5
$ws.Cells.Item(2, 13).Value = 0.102926
$ws.Cells.Item(2, 14).Value = 0.205852
$ws.Cells.Item(2, 15).Value = 0.03821605677617539
$ws.Cells.Item(2, 16).Value = 0.02724178813968781
$ws.Cells.Item(2, 17).Value = 3.257142880332
$ws.Cells.Item(2, 18).Value = 13.028571521328
$ws.Cells.Item(2, 19).Value = 0.001358955613465995
$ws.Cells.Item(2, 20).Value = 0.0006670046705825946

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Thbs1"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 31.645482
$ws.Cells.Item(3, 8).Value = 63.290964
$ws.Cells.Item(3, 9).Value = 0.03555980726701226
$ws.Cells.Item(3, 10).Value = 0.0244846141215985
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.169947333333333
$ws.Cells.Item(3, 14).Value = 6.509842000000001
$ws.Cells.Item(3, 15).Value = 0.8056937070514454
$ws.Cells.Item(3, 16).Value = 0.8614914433031574
$ws.Cells.Item(3, 17).Value = 68.66902927794801
$ws.Cells.Item(3, 18).Value = 412.0141756676881
$ws.Cells.Item(3, 19).Value = 0.02865031293899403
$ws.Cells.Item(3, 20).Value = 0.02109328555833676

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Thbs1"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 31.645482
$ws.Cells.Item(4, 8).Value = 63.290964
$ws.Cells.Item(4, 9).Value = 0.03555980726701226
$ws.Cells.Item(4, 10).Value = 0.0244846141215985
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.4203924999999999
$ws.Cells.Item(4, 14).Value = 0.8407849999999999
$ws.Cells.Item(4, 15).Value = 0.1560902361723793
$ws.Cells.Item(4, 16).Value = 0.1112667685571547
$ws.Cells.Item(4, 17).Value = 13.303523291685
$ws.Cells.Item(4, 18).Value = 53.21409316673999
$ws.Cells.Item(4, 19).Value = 0.005550538714552233
$ws.Cells.Item(4, 20).Value = 0.002724323892679142

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Thbs1"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 83.22744366666667
$ws.Cells.Item(5, 8).Value = 249.682331
$ws.Cells.Item(5, 9).Value = 0.09352209759714789
$ws.Cells.Item(5, 10).Value = 0.09659160077758068
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.5
$ws.Cells.Item(5, 13).Value = 0.102926
$ws.Cells.Item(5, 14).Value = 0.205852
$ws.Cells.Item(5, 15).Value = 0.03821605677617539
$ws.Cells.Item(5, 16).Value = 0.02724178813968781
$ws.Cells.Item(5, 17).Value = 8.566267866835334
$ws.Cells.Item(5, 18).Value = 51.397607201012
$ws.Cells.Item(5, 19).Value = 0.00357404579159962
$ws.Cells.Item(5, 20).Value = 0.002631327924456157

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Thbs1"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 83.22744366666667
$ws.Cells.Item(6, 8).Value = 249.682331
$ws.Cells.Item(6, 9).Value = 0.09352209759714789
$ws.Cells.Item(6, 10).Value = 0.09659160077758068
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.169947333333333
$ws.Cells.Item(6, 14).Value = 6.509842000000001
$ws.Cells.Item(6, 15).Value = 0.8056937070514454
$ws.Cells.Item(6, 16).Value = 0.8614914433031574
$ws.Cells.Item(6, 17).Value = 180.5991694446336
$ws.Cells.Item(6, 18).Value = 1625.392525001702
$ws.Cells.Item(6, 19).Value = 0.07535016550427315
$ws.Cells.Item(6, 20).Value = 0.08321283756484037

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Thbs1"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 83.22744366666667
$ws.Cells.Item(7, 8).Value = 249.682331
$ws.Cells.Item(7, 9).Value = 0.09352209759714789
$ws.Cells.Item(7, 10).Value = 0.09659160077758068
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.4203924999999999
$ws.Cells.Item(7, 14).Value = 0.8407849999999999
$ws.Cells.Item(7, 15).Value = 0.1560902361723793
$ws.Cells.Item(7, 16).Value = 0.1112667685571547
$ws.Cells.Item(7, 17).Value = 34.98819311163916
$ws.Cells.Item(7, 18).Value = 209.929158669835
$ws.Cells.Item(7, 19).Value = 0.01459788630127512
$ws.Cells.Item(7, 20).Value = 0.01074743528828415

# Row 8
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Thbs1"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 524.5768889999999
$ws.Cells.Item(8, 8).Value = 1573.730667
$ws.Cells.Item(8, 9).Value = 0.5894633891046084
$ws.Cells.Item(8, 10).Value = 0.6088102578564109
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.5
$ws.Cells.Item(8, 13).Value = 0.102926
$ws.Cells.Item(8, 14).Value = 0.205852
$ws.Cells.Item(8, 15).Value = 0.03821605677617539
$ws.Cells.Item(8, 16).Value = 0.02724178813968781
$ws.Cells.Item(8, 17).Value = 53.992600877214
$ws.Cells.Item(8, 18).Value = 323.955605263284
$ws.Cells.Item(8, 19).Value = 0.02252696634549848
$ws.Cells.Item(8, 20).Value = 0.01658508006179305

# Row 9
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Thbs1"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 524.5768889999999
$ws.Cells.Item(9, 8).Value = 1573.730667
$ws.Cells.Item(9, 9).Value = 0.5894633891046084
$ws.Cells.Item(9, 10).Value = 0.6088102578564109
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.169947333333333
$ws.Cells.Item(9, 14).Value = 6.509842000000001
$ws.Cells.Item(9, 15).Value = 0.8056937070514454
$ws.Cells.Item(9, 16).Value = 0.8614914433031574
$ws.Cells.Item(9, 17).Value = 1138.304221413846
$ws.Cells.Item(9, 18).Value = 10244.73799272461
$ws.Cells.Item(9, 19).Value = 0.4749269431388005
$ws.Cells.Item(9, 20).Value = 0.5244848277384868

# Row 10
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Thbs1"
$ws.Cells.Item(10, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 524.5768889999999
$ws.Cells.Item(10, 8).Value = 1573.730667
$ws.Cells.Item(10, 9).Value = 0.5894633891046084
$ws.Cells.Item(10, 10).Value = 0.6088102578564109
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.4203924999999999
$ws.Cells.Item(10, 14).Value = 0.8407849999999999
$ws.Cells.Item(10, 15).Value = 0.1560902361723793
$ws.Cells.Item(10, 16).Value = 0.1112667685571547
$ws.Cells.Item(10, 17).Value = 220.5281898089324
$ws.Cells.Item(10, 18).Value = 1323.169138853595
$ws.Cells.Item(10, 19).Value = 0.09200947962030943
$ws.Cells.Item(10, 20).Value = 0.06774035005613094

# Row 11
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Thbs1"
$ws.Cells.Item(11, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 181.4813383333333
$ws.Cells.Item(11, 8).Value = 544.4440149999999
$ws.Cells.Item(11, 9).Value = 0.2039293133121744
$ws.Cells.Item(11, 10).Value = 0.210622508737405
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.5
$ws.Cells.Item(11, 13).Value = 0.102926
$ws.Cells.Item(11, 14).Value = 0.205852
$ws.Cells.Item(11, 15).Value = 0.03821605677617539
$ws.Cells.Item(11, 16).Value = 0.02724178813968781
$ws.Cells.Item(11, 17).Value = 18.67914822929666
$ws.Cells.Item(11, 18).Value = 112.07488937578
$ws.Cells.Item(11, 19).Value = 0.007793374215864517
$ws.Cells.Item(11, 20).Value = 0.00573773376047393

# Row 12
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Thbs1"
$ws.Cells.Item(12, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 181.4813383333333
$ws.Cells.Item(12, 8).Value = 544.4440149999999
$ws.Cells.Item(12, 9).Value = 0.2039293133121744
$ws.Cells.Item(12, 10).Value = 0.210622508737405
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 2.169947333333333
$ws.Cells.Item(12, 14).Value = 6.509842000000001
$ws.Cells.Item(12, 15).Value = 0.8056937070514454
$ws.Cells.Item(12, 16).Value = 0.8614914433031574
$ws.Cells.Item(12, 17).Value = 393.8049461661811
$ws.Cells.Item(12, 18).Value = 3544.24451549563
$ws.Cells.Item(12, 19).Value = 0.1643045644189415
$ws.Cells.Item(12, 20).Value = 0.1814494890443189

# Row 13
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Thbs1"
$ws.Cells.Item(13, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 181.4813383333333
$ws.Cells.Item(13, 8).Value = 544.4440149999999
$ws.Cells.Item(13, 9).Value = 0.2039293133121744
$ws.Cells.Item(13, 10).Value = 0.210622508737405
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.4203924999999999
$ws.Cells.Item(13, 14).Value = 0.8407849999999999
$ws.Cells.Item(13, 15).Value = 0.1560902361723793
$ws.Cells.Item(13, 16).Value = 0.1112667685571547
$ws.Cells.Item(13, 17).Value = 76.29339352529581
$ws.Cells.Item(13, 18).Value = 457.7603611517749
$ws.Cells.Item(13, 19).Value = 0.03183137467736843
$ws.Cells.Item(13, 20).Value = 0.02343528593261213

# Row 14
$ws.Cells.Item(14, 1).Value = "Neutro"
$ws.Cells.Item(14, 2).Value = "Thbs1"
$ws.Cells.Item(14, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 15.79677433333333
$ws.Cells.Item(14, 8).Value = 47.390323
$ws.Cells.Item(14, 9).Value = 0.01775072507139627
$ws.Cells.Item(14, 10).Value = 0.01833332435500452
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.5
$ws.Cells.Item(14, 13).Value = 0.102926
$ws.Cells.Item(14, 14).Value = 0.205852
$ws.Cells.Item(14, 15).Value = 0.03821605677617539
$ws.Cells.Item(14, 16).Value = 0.02724178813968781
$ws.Cells.Item(14, 17).Value = 1.625898795032667
$ws.Cells.Item(14, 18).Value = 9.755392770196
$ws.Cells.Item(14, 19).Value = 0.0006783627171467598
$ws.Cells.Item(14, 20).Value = 0.0004994325379752116

# Row 15
$ws.Cells.Item(15, 1).Value = "Neutro"
$ws.Cells.Item(15, 2).Value = "Thbs1"
$ws.Cells.Item(15, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 15.79677433333333
$ws.Cells.Item(15, 8).Value = 47.390323
$ws.Cells.Item(15, 9).Value = 0.01775072507139627
$ws.Cells.Item(15, 10).Value = 0.01833332435500452
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2.169947333333333
$ws.Cells.Item(15, 14).Value = 6.509842000000001
$ws.Cells.Item(15, 15).Value = 0.8056937070514454
$ws.Cells.Item(15, 16).Value = 0.8614914433031574
$ws.Cells.Item(15, 17).Value = 34.27816833988511
$ws.Cells.Item(15, 18).Value = 308.5035150589661
$ws.Cells.Item(15, 19).Value = 0.01430164748562429
$ws.Cells.Item(15, 20).Value = 0.01579400205913777

# Row 16
$ws.Cells.Item(16, 1).Value = "Neutro"
$ws.Cells.Item(16, 2).Value = "Thbs1"
$ws.Cells.Item(16, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 15.79677433333333
$ws.Cells.Item(16, 8).Value = 47.390323
$ws.Cells.Item(16, 9).Value = 0.01775072507139627
$ws.Cells.Item(16, 10).Value = 0.01833332435500452
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.4203924999999999
$ws.Cells.Item(16, 14).Value = 0.8407849999999999
$ws.Cells.Item(16, 15).Value = 0.1560902361723793
$ws.Cells.Item(16, 16).Value = 0.1112667685571547
$ws.Cells.Item(16, 17).Value = 6.640845453925833
$ws.Cells.Item(16, 18).Value = 39.845072723555
$ws.Cells.Item(16, 19).Value = 0.002770714868625218
$ws.Cells.Item(16, 20).Value = 0.002039889757891535

# Row 17
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Thbs1"
$ws.Cells.Item(17, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 2
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 53.1948375
$ws.Cells.Item(17, 8).Value = 106.389675
$ws.Cells.Item(17, 9).Value = 0.05977466764766092
$ws.Cells.Item(17, 10).Value = 0.0411576941520005
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.5
$ws.Cells.Item(17, 13).Value = 0.102926
$ws.Cells.Item(17, 14).Value = 0.205852
$ws.Cells.Item(17, 15).Value = 0.03821605677617539
$ws.Cells.Item(17, 16).Value = 0.02724178813968781
$ws.Cells.Item(17, 17).Value = 5.475131844525
$ws.Cells.Item(17, 18).Value = 21.9005273781
$ws.Cells.Item(17, 19).Value = 0.002284352092600024
$ws.Cells.Item(17, 20).Value = 0.001121209184406865

# Row 18
$ws.Cells.Item(18, 1).Value = "sCs"
$ws.Cells.Item(18, 2).Value = "Thbs1"
$ws.Cells.Item(18, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 2
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 53.1948375
$ws.Cells.Item(18, 8).Value = 106.389675
$ws.Cells.Item(18, 9).Value = 0.05977466764766092
$ws.Cells.Item(18, 10).Value = 0.0411576941520005
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 2.169947333333333
$ws.Cells.Item(18, 14).Value = 6.509842000000001
$ws.Cells.Item(18, 15).Value = 0.8056937070514454
$ws.Cells.Item(18, 16).Value = 0.8614914433031574
$ws.Cells.Item(18, 17).Value = 115.429995780225
$ws.Cells.Item(18, 18).Value = 692.5799746813501
$ws.Cells.Item(18, 19).Value = 0.04816007356481203
$ws.Cells.Item(18, 20).Value = 0.03545700133803684

# Row 19
$ws.Cells.Item(19, 1).Value = "sCs"
$ws.Cells.Item(19, 2).Value = "Thbs1"
$ws.Cells.Item(19, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 2
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 53.1948375
$ws.Cells.Item(19, 8).Value = 106.389675
$ws.Cells.Item(19, 9).Value = 0.05977466764766092
$ws.Cells.Item(19, 10).Value = 0.0411576941520005
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.4203924999999999
$ws.Cells.Item(19, 14).Value = 0.8407849999999999
$ws.Cells.Item(19, 15).Value = 0.1560902361723793
$ws.Cells.Item(19, 16).Value = 0.1112667685571547
$ws.Cells.Item(19, 17).Value = 22.36271072371875
$ws.Cells.Item(19, 18).Value = 89.45084289487498
$ws.Cells.Item(19, 19).Value = 0.009330241990248874
$ws.Cells.Item(19, 20).Value = 0.0045794836295568
